# Applies the Unicorn_Profits price-refresh update (scheduled runner data sync).
# For each affected Leve row, writes the refreshed currentAveragePrice* / LevePrice* /
# LeveProfit* columns (H-N). Cells that the refresh no longer populates are cleared
# (ClearContents) so the saved XML omits the <c> element entirely, matching upstream.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 137
$ws.Range("H137").Value = 2561892
$ws.Range("I137").Value = 3713476.5
$ws.Range("J137").Value = 1660651.9
$ws.Range("K137").Value = 11140429.5
$ws.Range("L137").Value = 4981955.699999999
$ws.Range("M137").Value = -11137879.5
$ws.Range("N137").Value = -4987055.699999999

# Row 138
$ws.Range("H138").Value = 2147.0698
$ws.Range("I138").Value = 1332.963
$ws.Range("J138").Value = 3520.875
$ws.Range("K138").Value = 3998.889
$ws.Range("L138").Value = 10562.625
$ws.Range("M138").Value = 1141.111
$ws.Range("N138").Value = -20842.625

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 1441007.9
$ws.Range("I32").Value = 1441007.9
$ws.Range("K32").Value = 1441007.9
$ws.Range("M32").Value = -1440720.9

# Row 74
$ws.Range("H74").Value = 15403257
$ws.Range("I74").Value = 13734328
$ws.Range("J74").Value = 17599216
$ws.Range("K74").Value = 13734328
$ws.Range("L74").Value = 17599216
$ws.Range("M74").Value = -13733454
$ws.Range("N74").Value = -17600964

# Row 77
$ws.Range("H77").Value = 15403257
$ws.Range("I77").Value = 13734328
$ws.Range("J77").Value = 17599216
$ws.Range("K77").Value = 68671640
$ws.Range("L77").Value = 87996080
$ws.Range("M77").Value = -68667272
$ws.Range("N77").Value = -88004816

$ws = $wb.Worksheets.Item("CRP")
# Row 6
$ws.Range("H6").Value = 67677000
$ws.Range("J6").Value = 100515500
$ws.Range("L6").Value = 100515500
$ws.Range("N6").Value = -100515726

# Row 16
$ws.Range("H16").Value = 889.55
$ws.Range("I16").Value = 786.63635
$ws.Range("J16").Value = 1015.3333
$ws.Range("K16").Value = 786.63635
$ws.Range("L16").Value = 1015.3333
$ws.Range("M16").Value = -499.63635
$ws.Range("N16").Value = -1589.3333

# Row 17
$ws.Range("H17").Value = 3000
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()

# Row 25
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents()

# Row 31
$ws.Range("H31").Value = 1785909.1
$ws.Range("I31").Value = 1191.8462
$ws.Range("J31").Value = 2396470.2
$ws.Range("K31").Value = 1191.8462
$ws.Range("L31").Value = 2396470.2
$ws.Range("M31").Value = -896.8462
$ws.Range("N31").Value = -2397060.2

# Row 34
$ws.Range("H34").Value = 1785909.1
$ws.Range("I34").Value = 1191.8462
$ws.Range("J34").Value = 2396470.2
$ws.Range("K34").Value = 1191.8462
$ws.Range("L34").Value = 2396470.2
$ws.Range("M34").Value = -989.8462
$ws.Range("N34").Value = -2396874.2

# Row 41
$ws.Range("H41").Value = 8200
$ws.Range("J41").Value = 15000
$ws.Range("L41").Value = 15000
$ws.Range("N41").Value = -15856

# Row 58
$ws.Range("H58").Value = 3087.0984
$ws.Range("I58").Value = 3285.4866
$ws.Range("J58").Value = 2781.25
$ws.Range("K58").Value = 3285.4866
$ws.Range("L58").Value = 2781.25
$ws.Range("M58").Value = -3082.4866
$ws.Range("N58").Value = -3187.25

# Row 86
$ws.Range("H86").Value = 3355.394
$ws.Range("I86").Value = 2144
$ws.Range("J86").Value = 5219.077
$ws.Range("K86").Value = 2144
$ws.Range("L86").Value = 5219.077
$ws.Range("M86").Value = -1021
$ws.Range("N86").Value = -7465.077

# Row 89
$ws.Range("H89").Value = 3355.394
$ws.Range("I89").Value = 2144
$ws.Range("J89").Value = 5219.077
$ws.Range("K89").Value = 10720
$ws.Range("L89").Value = 26095.385
$ws.Range("M89").Value = -5104
$ws.Range("N89").Value = -37327.385

# Row 92
$ws.Range("H92").Value = 20200.334
$ws.Range("J92").Value = 20200.334
$ws.Range("L92").Value = 20200.334
$ws.Range("N92").Value = -25192.334

# Row 113
$ws.Range("H113").Value = 889.55
$ws.Range("I113").Value = 786.63635
$ws.Range("J113").Value = 1015.3333
$ws.Range("K113").Value = 786.63635
$ws.Range("L113").Value = 1015.3333
$ws.Range("M113").Value = 1383.36365
$ws.Range("N113").Value = -5355.3333

# Row 136
$ws.Range("H136").Value = 3087.0984
$ws.Range("I136").Value = 3285.4866
$ws.Range("J136").Value = 2781.25
$ws.Range("K136").Value = 9856.459800000001
$ws.Range("L136").Value = 8343.75
$ws.Range("M136").Value = -7306.459800000001
$ws.Range("N136").Value = -13443.75

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 800.09753
$ws.Range("I5").Value = 486.4
$ws.Range("J5").Value = 1290.25
$ws.Range("K5").Value = 1459.2
$ws.Range("L5").Value = 3870.75
$ws.Range("M5").Value = -1347.2
$ws.Range("N5").Value = -4094.75

# Row 9
$ws.Range("H9").Value = 78750440
$ws.Range("I9").Value = 901
$ws.Range("J9").Value = 105000290
$ws.Range("K9").Value = 2703
$ws.Range("L9").Value = 315000870
$ws.Range("M9").Value = -2479
$ws.Range("N9").Value = -315001318

# Row 58
$ws.Range("H58").Value = 2539.125
$ws.Range("J58").Value = 3222.6
$ws.Range("L58").Value = 9667.799999999999
$ws.Range("N58").Value = -9923.799999999999

# Row 68
$ws.Range("H68").Value = 862.3875
$ws.Range("I68").Value = 510.27585
$ws.Range("J68").Value = 1062.6078
$ws.Range("K68").Value = 1530.82755
$ws.Range("L68").Value = 3187.8234
$ws.Range("M68").Value = -719.82755
$ws.Range("N68").Value = -4809.8234

# Row 71
$ws.Range("H71").Value = 862.3875
$ws.Range("I71").Value = 510.27585
$ws.Range("J71").Value = 1062.6078
$ws.Range("K71").Value = 4592.48265
$ws.Range("L71").Value = 9563.4702
$ws.Range("M71").Value = -536.4826499999999
$ws.Range("N71").Value = -17675.4702

# Row 76
$ws.Range("H76").Value = 6062.5
$ws.Range("J76").Value = 6062.5
$ws.Range("L76").Value = 18187.5
$ws.Range("N76").Value = -18953.5

# Row 79
$ws.Range("H79").Value = 6062.5
$ws.Range("J79").Value = 6062.5
$ws.Range("L79").Value = 18187.5
$ws.Range("N79").Value = -20839.5

# Row 107
$ws.Range("H107").Value = 1027.6216
$ws.Range("I107").Value = 532.44446
$ws.Range("J107").Value = 1186.7858
$ws.Range("K107").Value = 1597.33338
$ws.Range("L107").Value = 3560.3574
$ws.Range("M107").Value = 322.66662
$ws.Range("N107").Value = -7400.357400000001

# Row 131
$ws.Range("H131").Value = 1062.8534
$ws.Range("I131").Value = 743.6667
$ws.Range("J131").Value = 1163.6492
$ws.Range("K131").Value = 2231.0001
$ws.Range("L131").Value = 3490.9476
$ws.Range("M131").Value = 2808.9999
$ws.Range("N131").Value = -13570.9476

# Row 135
$ws.Range("H135").Value = 800.09753
$ws.Range("I135").Value = 486.4
$ws.Range("J135").Value = 1290.25
$ws.Range("K135").Value = 4377.599999999999
$ws.Range("L135").Value = 11612.25
$ws.Range("M135").Value = -1842.599999999999
$ws.Range("N135").Value = -16682.25

$ws = $wb.Worksheets.Item("GSM")
# Row 12
$ws.Range("H12").Value = 6696.4614
$ws.Range("I12").Value = 1420.8334
$ws.Range("J12").Value = 70004
$ws.Range("K12").Value = 1420.8334
$ws.Range("L12").Value = 70004
$ws.Range("M12").Value = -1280.8334
$ws.Range("N12").Value = -70284

# Row 42
$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()

# Row 109
$ws.Range("H109").Value = 10285
$ws.Range("J109").Value = 10285
$ws.Range("L109").Value = 10285
$ws.Range("N109").Value = -12365

# Row 115
$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# Row 15
$ws.Range("H15").Value = 4000
$ws.Range("J15").Value = 4000
$ws.Range("L15").Value = 4000
$ws.Range("N15").Value = -4576

# Row 122
$ws.Range("H122").Value = 57406.445
$ws.Range("I122").Value = 113079
$ws.Range("J122").Value = 1733.8889
$ws.Range("K122").Value = 339237
$ws.Range("L122").Value = 5201.6667
$ws.Range("M122").Value = -336787
$ws.Range("N122").Value = -10101.6667

Write-Output "Applied all Unicorn_Profits updates"